$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 4347916
$ws.Range("I53").Value = 112.5
$ws.Range("J53").Value = 9090974
$ws.Range("K53").Value = 112.5
$ws.Range("L53").Value = 9090974
$ws.Range("M53").Value = 524.5
$ws.Range("N53").Value = -9092248

$ws.Range("H129").Value = 1179.3914
$ws.Range("I129").Value = 660.1818
$ws.Range("J129").Value = 1342.5714
$ws.Range("K129").Value = 1980.5454
$ws.Range("L129").Value = 4027.7142
$ws.Range("M129").Value = 3019.4546
$ws.Range("N129").Value = -14027.7142

$ws.Range("H132").Value = 2047.7812
$ws.Range("I132").Value = 1368.762
$ws.Range("K132").Value = 4106.286
$ws.Range("M132").Value = -1576.286

$ws.Range("H138").Value = 3278.4
$ws.Range("I138").Value = 2532.4827
$ws.Range("J138").Value = 6883.6665
$ws.Range("K138").Value = 7597.4481
$ws.Range("L138").Value = 20650.9995
$ws.Range("M138").Value = -2457.4481
$ws.Range("N138").Value = -30930.9995

$ws.Range("H141").Value = 2907.02
$ws.Range("I141").Value = 1659.1945
$ws.Range("K141").Value = 4977.583500000001
$ws.Range("M141").Value = 202.4164999999994


# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 41666976
$ws.Range("I5").Value = 62500090
$ws.Range("J5").Value = 750
$ws.Range("K5").Value = 62500090
$ws.Range("L5").Value = 750
$ws.Range("M5").Value = -62499978
$ws.Range("N5").Value = -974

$ws.Range("H11").Value = 4201.5
$ws.Range("I11").Value = 2403
$ws.Range("J11").Value = 6000
$ws.Range("K11").Value = 2403
$ws.Range("L11").Value = 6000
$ws.Range("M11").Value = -2259
$ws.Range("N11").Value = -6288

$ws.Range("H19").Value = 6236.7144
$ws.Range("I19").Value = 1778.5
$ws.Range("J19").Value = 8020
$ws.Range("K19").Value = 1778.5
$ws.Range("L19").Value = 8020
$ws.Range("M19").Value = -1549.5
$ws.Range("N19").Value = -8478

$ws.Range("H32").Value = 6870.0396
$ws.Range("I32").Value = 2876.3386
$ws.Range("J32").Value = 24556.428
$ws.Range("K32").Value = 2876.3386
$ws.Range("L32").Value = 24556.428
$ws.Range("M32").Value = -2589.3386
$ws.Range("N32").Value = -25130.428

$ws.Range("H61").Value = 1837.625
$ws.Range("J61").Value = 2095.875
$ws.Range("L61").Value = 2095.875
$ws.Range("N61").Value = -2519.875

$ws.Range("H88").Value = 1440
$ws.Range("I88").Value = 1433.3334
$ws.Range("J88").Value = 1445
$ws.Range("K88").Value = 1433.3334
$ws.Range("L88").Value = 1445
$ws.Range("M88").Value = -1027.3334
$ws.Range("N88").Value = -2257

$ws.Range("H91").Value = 1440
$ws.Range("I91").Value = 1433.3334
$ws.Range("J91").Value = 1445
$ws.Range("K91").Value = 1433.3334
$ws.Range("L91").Value = 1445
$ws.Range("M91").Value = -29.33339999999998
$ws.Range("N91").Value = -4253

$ws.Range("H136").Value = 1837.625
$ws.Range("J136").Value = 2095.875
$ws.Range("L136").Value = 6287.625
$ws.Range("N136").Value = -11387.625


# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 41666976
$ws.Range("I4").Value = 62500090
$ws.Range("J4").Value = 750
$ws.Range("K4").Value = 62500090
$ws.Range("L4").Value = 750
$ws.Range("M4").Value = -62499975
$ws.Range("N4").Value = -980

$ws.Range("H86").Value = 1801.6666
$ws.Range("I86").Value = 1735.1666
$ws.Range("J86").Value = 2001.1666
$ws.Range("K86").Value = 1735.1666
$ws.Range("L86").Value = 2001.1666
$ws.Range("M86").Value = -612.1666
$ws.Range("N86").Value = -4247.1666

$ws.Range("H89").Value = 1801.6666
$ws.Range("I89").Value = 1735.1666
$ws.Range("J89").Value = 2001.1666
$ws.Range("K89").Value = 8675.833000000001
$ws.Range("L89").Value = 10005.833
$ws.Range("M89").Value = -3059.833000000001
$ws.Range("N89").Value = -21237.833


# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1669.7258
$ws.Range("I31").Value = 1230.3572
$ws.Range("J31").Value = 2031.5588
$ws.Range("K31").Value = 1230.3572
$ws.Range("L31").Value = 2031.5588
$ws.Range("M31").Value = -935.3571999999999
$ws.Range("N31").Value = -2621.5588

$ws.Range("H34").Value = 1669.7258
$ws.Range("I34").Value = 1230.3572
$ws.Range("J34").Value = 2031.5588
$ws.Range("K34").Value = 1230.3572
$ws.Range("L34").Value = 2031.5588
$ws.Range("M34").Value = -1028.3572
$ws.Range("N34").Value = -2435.5588

$ws.Range("H59").Value = 53666.668
$ws.Range("I59").Value = 50000
$ws.Range("J59").Value = 55500
$ws.Range("K59").Value = 50000
$ws.Range("L59").Value = 55500
$ws.Range("M59").Value = -48855
$ws.Range("N59").Value = -57790

$ws.Range("H132").Value = 3390.68
$ws.Range("I132").Value = 2934.158
$ws.Range("K132").Value = 8802.474
$ws.Range("M132").Value = -6272.474


# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 981.3099999999999
$ws.Range("I68").Value = 772.5185
$ws.Range("J68").Value = 1226.4131
$ws.Range("K68").Value = 2317.5555
$ws.Range("L68").Value = 3679.2393
$ws.Range("M68").Value = -1506.5555
$ws.Range("N68").Value = -5301.2393

$ws.Range("H71").Value = 981.3099999999999
$ws.Range("I71").Value = 772.5185
$ws.Range("J71").Value = 1226.4131
$ws.Range("K71").Value = 6952.6665
$ws.Range("L71").Value = 11037.7179
$ws.Range("M71").Value = -2896.6665
$ws.Range("N71").Value = -19149.7179

$ws.Range("H102").Value = 3755.7144
$ws.Range("I102").Value = 2800
$ws.Range("J102").Value = 3915
$ws.Range("K102").Value = 8400
$ws.Range("L102").Value = 11745
$ws.Range("M102").Value = -5966
$ws.Range("N102").Value = -16613

$ws.Range("H113").Value = 9524195
$ws.Range("I113").Value = 427.66666
$ws.Range("J113").Value = 11905136
$ws.Range("K113").Value = 1282.99998
$ws.Range("L113").Value = 35715408
$ws.Range("M113").Value = 887.0000199999999
$ws.Range("N113").Value = -35719748

$ws.Range("H123").Value = 1450
$ws.Range("I123").Value = 1000
$ws.Range("J123").Value = 1900
$ws.Range("K123").Value = 3000
$ws.Range("L123").Value = 5700
$ws.Range("M123").Value = -550
$ws.Range("N123").Value = -10600


# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 5745
$ws.Range("J33").Value = 5745
$ws.Range("L33").Value = 5745
$ws.Range("N33").Value = -6249

$ws.Range("H70").Value = 5367.931
$ws.Range("I70").Value = 5030.769
$ws.Range("J70").Value = 5641.875
$ws.Range("K70").Value = 5030.769
$ws.Range("L70").Value = 5641.875
$ws.Range("M70").Value = -4760.769
$ws.Range("N70").Value = -6181.875

$ws.Range("H73").Value = 5367.931
$ws.Range("I73").Value = 5030.769
$ws.Range("J73").Value = 5641.875
$ws.Range("K73").Value = 5030.769
$ws.Range("L73").Value = 5641.875
$ws.Range("M73").Value = -4094.769
$ws.Range("N73").Value = -7513.875


# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1516.25
$ws.Range("I16").Value = 1743.4445
$ws.Range("J16").Value = 834.6667
$ws.Range("K16").Value = 834.6667
$ws.Range("L16").Value = 834.6667
$ws.Range("M16").Value = -1573.4445
$ws.Range("N16").Value = -1174.6667

$ws.Range("H46").Value = 964.4783
$ws.Range("I46").Value = 984.5
$ws.Range("J46").Value = 933.3333
$ws.Range("K46").Value = 984.5
$ws.Range("L46").Value = 933.3333
$ws.Range("M46").Value = -796.5
$ws.Range("N46").Value = -1309.3333

